$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.433944
$ws.Range("H2").Value = 4.301832
$ws.Range("I2").Value = 0.003882998715548277
$ws.Range("J2").Value = 0.003886188105009087
$ws.Range("M2").Value = 2.906846333333333
$ws.Range("N2").Value = 8.720538999999999
$ws.Range("O2").Value = 0.005520525738044089
$ws.Range("P2").Value = 0.005624540846623205
$ws.Range("Q2").Value = 4.168254858605333
$ws.Range("R2").Value = 37.51429372744799
$ws.Range("S2").Value = 0.0000214361943499764
$ws.Range("T2").Value = 0.00002185802373428484
$ws.Range("G3").Value = 1.433944
$ws.Range("H3").Value = 4.301832
$ws.Range("I3").Value = 0.003882998715548277
$ws.Range("J3").Value = 0.003886188105009087
$ws.Range("O3").Value = 0.3528665483720876
$ws.Range("P3").Value = 0.3595150912979765
$ws.Range("Q3").Value = 266.4307304203173
$ws.Range("R3").Value = 2397.876573782856
$ws.Range("S3").Value = 0.00137018035408877
$ws.Range("T3").Value = 0.001397143271373452
$ws.Range("G4").Value = 1.433944
$ws.Range("H4").Value = 4.301832
$ws.Range("I4").Value = 0.003882998715548277
$ws.Range("J4").Value = 0.003886188105009087
$ws.Range("M4").Value = 137.0717086666666
$ws.Range("N4").Value = 411.2151259999999
$ws.Range("O4").Value = 0.2603191943704447
$ws.Range("P4").Value = 0.2652240042658267
$ws.Range("Q4").Value = 196.5531542123146
$ws.Range("R4").Value = 1768.978387910832
$ws.Range("S4").Value = 0.001010819097372999
$ws.Range("T4").Value = 0.001030710370540735
$ws.Range("G5").Value = 1.433944
$ws.Range("H5").Value = 4.301832
$ws.Range("I5").Value = 0.003882998715548277
$ws.Range("J5").Value = 0.003886188105009087
$ws.Range("M5").Value = 29.2127365
$ws.Range("N5").Value = 58.425473
$ws.Range("O5").Value = 0.05547925319534149
$ws.Range("P5").Value = 0.03768304451958546
$ws.Range("Q5").Value = 41.889428227756
$ws.Range("R5").Value = 251.336569366536
$ws.Range("S5").Value = 0.0002154258688970886
$ws.Range("T5").Value = 0.0001464433993725409
$ws.Range("G6").Value = 1.433944
$ws.Range("H6").Value = 4.301832
$ws.Range("I6").Value = 0.003882998715548277
$ws.Range("J6").Value = 0.003886188105009087
$ws.Range("M6").Value = 171.5584106666666
$ws.Range("N6").Value = 514.6752319999999
$ws.Range("O6").Value = 0.3258144783240821
$ws.Range("P6").Value = 0.331953319069988
$ws.Range("Q6").Value = 246.0051536250027
$ws.Range("R6").Value = 2214.046382625024
$ws.Range("S6").Value = 0.001265137200839443
$ws.Range("T6").Value = 0.001290033039988074
$ws.Range("I7").Value = 0.7877262822264709
$ws.Range("J7").Value = 0.7883732991550308
$ws.Range("M7").Value = 2.906846333333333
$ws.Range("N7").Value = 8.720538999999999
$ws.Range("O7").Value = 0.005520525738044089
$ws.Range("P7").Value = 0.005624540846623205
$ws.Range("Q7").Value = 845.5948980858683
$ws.Range("R7").Value = 7610.354082772815
$ws.Range("S7").Value = 0.004348663215565015
$ws.Range("T7").Value = 0.004434237823484567
$ws.Range("I8").Value = 0.7877262822264709
$ws.Range("J8").Value = 0.7883732991550308
$ws.Range("O8").Value = 0.3528665483720876
$ws.Range("P8").Value = 0.3595150912979765
$ws.Range("S8").Value = 0.2779622542712318
$ws.Range("T8").Value = 0.2834320986226079
$ws.Range("I9").Value = 0.7877262822264709
$ws.Range("J9").Value = 0.7883732991550308
$ws.Range("M9").Value = 137.0717086666666
$ws.Range("N9").Value = 411.2151259999999
$ws.Range("O9").Value = 0.2603191943704447
$ws.Range("P9").Value = 0.2652240042658267
$ws.Range("Q9").Value = 39873.84410084486
$ws.Range("R9").Value = 358864.5969076037
$ws.Range("S9").Value = 0.2050602711736204
$ws.Range("T9").Value = 0.2090955232581578
$ws.Range("I10").Value = 0.7877262822264709
$ws.Range("J10").Value = 0.7883732991550308
$ws.Range("M10").Value = 29.2127365
$ws.Range("N10").Value = 58.425473
$ws.Range("O10").Value = 0.05547925319534149
$ws.Range("P10").Value = 0.03768304451958546
$ws.Range("Q10").Value = 8497.917712492368
$ws.Range("R10").Value = 50987.50627495421
$ws.Range("S10").Value = 0.0437024658602674
$ws.Range("T10").Value = 0.02970830613011149
$ws.Range("I11").Value = 0.7877262822264709
$ws.Range("J11").Value = 0.7883732991550308
$ws.Range("M11").Value = 171.5584106666666
$ws.Range("N11").Value = 514.6752319999999
$ws.Range("O11").Value = 0.3258144783240821
$ws.Range("P11").Value = 0.331953319069988
$ws.Range("Q11").Value = 49905.94622079675
$ws.Range("R11").Value = 449153.5159871707
$ws.Range("S11").Value = 0.2566526277057863
$ws.Range("T11").Value = 0.2617031333206691
$ws.Range("G12").Value = 54.70735966666666
$ws.Range("H12").Value = 164.122079
$ws.Range("I12").Value = 0.1481428893434501
$ws.Range("J12").Value = 0.1482645698807303
$ws.Range("M12").Value = 2.906846333333333
$ws.Range("N12").Value = 8.720538999999999
$ws.Range("O12").Value = 0.005520525738044089
$ws.Range("P12").Value = 0.005624540846623205
$ws.Range("Q12").Value = 159.0258878533978
$ws.Range("R12").Value = 1431.232990680581
$ws.Range("S12").Value = 0.0008178266335287338
$ws.Range("T12").Value = 0.0008339201294011879
$ws.Range("G13").Value = 54.70735966666666
$ws.Range("H13").Value = 164.122079
$ws.Range("I13").Value = 0.1481428893434501
$ws.Range("J13").Value = 0.1482645698807303
$ws.Range("O13").Value = 0.3528665483720876
$ws.Range("P13").Value = 0.3595150912979765
$ws.Range("Q13").Value = 10164.77756129738
$ws.Range("R13").Value = 91482.9980516764
$ws.Range("S13").Value = 0.05227467002849137
$ws.Range("T13").Value = 0.05330335037692596
$ws.Range("G14").Value = 54.70735966666666
$ws.Range("H14").Value = 164.122079
$ws.Range("I14").Value = 0.1481428893434501
$ws.Range("J14").Value = 0.1482645698807303
$ws.Range("M14").Value = 137.0717086666666
$ws.Range("N14").Value = 411.2151259999999
$ws.Range("O14").Value = 0.2603191943704447
$ws.Range("P14").Value = 0.2652240042658267
$ws.Range("Q14").Value = 7498.831266151881
$ws.Range("R14").Value = 67489.48139536694
$ws.Range("S14").Value = 0.03856443760559686
$ws.Range("T14").Value = 0.03932332291451776
$ws.Range("G15").Value = 54.70735966666666
$ws.Range("H15").Value = 164.122079
$ws.Range("I15").Value = 0.1481428893434501
$ws.Range("J15").Value = 0.1482645698807303
$ws.Range("M15").Value = 29.2127365
$ws.Range("N15").Value = 58.425473
$ws.Range("O15").Value = 0.05547925319534149
$ws.Range("P15").Value = 0.03768304451958546
$ws.Range("Q15").Value = 1598.151682553061
$ws.Range("R15").Value = 9588.910095318366
$ws.Range("S15").Value = 0.008218856866974726
$ws.Range("T15").Value = 0.005587060387492747
$ws.Range("G16").Value = 54.70735966666666
$ws.Range("H16").Value = 164.122079
$ws.Range("I16").Value = 0.1481428893434501
$ws.Range("J16").Value = 0.1482645698807303
$ws.Range("M16").Value = 171.5584106666666
$ws.Range("N16").Value = 514.6752319999999
$ws.Range("O16").Value = 0.3258144783240821
$ws.Range("P16").Value = 0.331953319069988
$ws.Range("Q16").Value = 9385.507676183035
$ws.Range("R16").Value = 84469.56908564731
$ws.Range("S16").Value = 0.04826709820885842
$ws.Range("T16").Value = 0.04921691607239258
$ws.Range("G17").Value = 0.909222
$ws.Range("H17").Value = 1.818444
$ws.Range("I17").Value = 0.002462096049879378
$ws.Range("J17").Value = 0.001642745565709015
$ws.Range("M17").Value = 2.906846333333333
$ws.Range("N17").Value = 8.720538999999999
$ws.Range("O17").Value = 0.005520525738044089
$ws.Range("P17").Value = 0.005624540846623205
$ws.Range("Q17").Value = 2.642968636885999
$ws.Range("R17").Value = 15.857811821316
$ws.Range("S17").Value = 0.00001359206461289579
$ws.Range("T17").Value = 0.000009239689534939501
$ws.Range("G18").Value = 0.909222
$ws.Range("H18").Value = 1.818444
$ws.Range("I18").Value = 0.002462096049879378
$ws.Range("J18").Value = 0.001642745565709015
$ws.Range("O18").Value = 0.3528665483720876
$ws.Range("P18").Value = 0.3595150912979765
$ws.Range("Q18").Value = 168.935942808242
$ws.Range("R18").Value = 1013.615656849452
$ws.Range("S18").Value = 0.0008687913348814876
$ws.Range("T18").Value = 0.0005905918220352227
$ws.Range("G19").Value = 0.909222
$ws.Range("H19").Value = 1.818444
$ws.Range("I19").Value = 0.002462096049879378
$ws.Range("J19").Value = 0.001642745565709015
$ws.Range("M19").Value = 137.0717086666666
$ws.Range("N19").Value = 411.2151259999999
$ws.Range("O19").Value = 0.2603191943704447
$ws.Range("P19").Value = 0.2652240042658267
$ws.Range("Q19").Value = 124.628613097324
$ws.Range("R19").Value = 747.7716785839439
$ws.Range("S19").Value = 0.0006409308601672539
$ws.Range("T19").Value = 0.0004356955569272757
$ws.Range("G20").Value = 0.909222
$ws.Range("H20").Value = 1.818444
$ws.Range("I20").Value = 0.002462096049879378
$ws.Range("J20").Value = 0.001642745565709015
$ws.Range("M20").Value = 29.2127365
$ws.Range("N20").Value = 58.425473
$ws.Range("O20").Value = 0.05547925319534149
$ws.Range("P20").Value = 0.03768304451958546
$ws.Range("Q20").Value = 26.560862706003
$ws.Range("R20").Value = 106.243450824012
$ws.Range("S20").Value = 0.0001365952501425081
$ws.Range("T20").Value = 0.00006190365428696442
$ws.Range("G21").Value = 0.909222
$ws.Range("H21").Value = 1.818444
$ws.Range("I21").Value = 0.002462096049879378
$ws.Range("J21").Value = 0.001642745565709015
$ws.Range("M21").Value = 171.5584106666666
$ws.Range("N21").Value = 514.6752319999999
$ws.Range("O21").Value = 0.3258144783240821
$ws.Range("P21").Value = 0.331953319069988
$ws.Range("Q21").Value = 155.984681263168
$ws.Range("R21").Value = 935.9080875790079
$ws.Range("S21").Value = 0.0008021865400752329
$ws.Range("T21").Value = 0.0005453148429246127
$ws.Range("G22").Value = 21.33956566666667
$ws.Range("H22").Value = 64.018697
$ws.Range("I22").Value = 0.05778573366465133
$ws.Range("J22").Value = 0.05783319729352075
$ws.Range("M22").Value = 2.906846333333333
$ws.Range("N22").Value = 8.720538999999999
$ws.Range("O22").Value = 0.005520525738044089
$ws.Range("P22").Value = 0.005624540846623205
$ws.Range("Q22").Value = 62.03083821307588
$ws.Range("R22").Value = 558.2775439176829
$ws.Range("S22").Value = 0.0003190076299874684
$ws.Range("T22").Value = 0.0003252851804682261
$ws.Range("G23").Value = 21.33956566666667
$ws.Range("H23").Value = 64.018697
$ws.Range("I23").Value = 0.05778573366465133
$ws.Range("J23").Value = 0.05783319729352075
$ws.Range("O23").Value = 0.3528665483720876
$ws.Range("P23").Value = 0.3595150912979765
$ws.Range("Q23").Value = 3964.949863748045
$ws.Range("R23").Value = 35684.5487737324
$ws.Range("S23").Value = 0.02039065238339426
$ws.Range("T23").Value = 0.020791907205034
$ws.Range("G24").Value = 21.33956566666667
$ws.Range("H24").Value = 64.018697
$ws.Range("I24").Value = 0.05778573366465133
$ws.Range("J24").Value = 0.05783319729352075
$ws.Range("M24").Value = 137.0717086666666
$ws.Range("N24").Value = 411.2151259999999
$ws.Range("O24").Value = 0.2603191943704447
$ws.Range("P24").Value = 0.2652240042658267
$ws.Range("Q24").Value = 2925.050728134535
$ws.Range("R24").Value = 26325.45655321082
$ws.Range("S24").Value = 0.01504273563368711
$ws.Range("T24").Value = 0.01533875216568314
$ws.Range("G25").Value = 21.33956566666667
$ws.Range("H25").Value = 64.018697
$ws.Range("I25").Value = 0.05778573366465133
$ws.Range("J25").Value = 0.05783319729352075
$ws.Range("M25").Value = 29.2127365
$ws.Range("N25").Value = 58.425473
$ws.Range("O25").Value = 0.05547925319534149
$ws.Range("P25").Value = 0.03768304451958546
$ws.Range("Q25").Value = 623.3871088447802
$ws.Range("R25").Value = 3740.322653068681
$ws.Range("S25").Value = 0.003205909349059759
$ws.Range("T25").Value = 0.002179330948321711
$ws.Range("G26").Value = 21.33956566666667
$ws.Range("H26").Value = 64.018697
$ws.Range("I26").Value = 0.05778573366465133
$ws.Range("J26").Value = 0.05783319729352075
$ws.Range("M26").Value = 171.5584106666666
$ws.Range("N26").Value = 514.6752319999999
$ws.Range("O26").Value = 0.3258144783240821
$ws.Range("P26").Value = 0.331953319069988
$ws.Range("Q26").Value = 3660.9819700903
$ws.Range("R26").Value = 32948.8377308127
$ws.Range("S26").Value = 0.01882742866852272
$ws.Range("T26").Value = 0.01919792179401366
